# LeanLargeMen_DLWeek2025.pptx — update the DQN backtest results table on
# slide 9: new Total Return / Sharpe Ratio figures, plus the table's new
# (built-in) style.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# Find the results table on the slide (don't assume a fixed shape index).
$tblShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tblShape = $shp
    }
}

$tbl = $tblShape.Table

# Table styles are applied, not assigned via the Style property.
$tbl.ApplyStyle("{5AD344C3-E686-44AB-9A52-C60CB1F23614}")

# Update the metric values by matching each row's label, so the edit is
# resilient to row ordering.
for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
    $label = $tbl.Cell($r, 1).Shape.TextFrame.TextRange.Text
    if ($label -eq "Total Return (%)") {
        $tbl.Cell($r, 2).Shape.TextFrame.TextRange.Text = "161.14"
    } elseif ($label -eq "Sharpe Ratio ") {
        $tbl.Cell($r, 2).Shape.TextFrame.TextRange.Text = "0.7025 "
    }
}
